# feat: develop color module
# Adds a new bronze-object record (id 21) to the quiz data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new record as row 23.
# Columns: id | name_chinese | name_pinyin | name_en | findspot | museum | start_time | end_time | image_path | url
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "错金嵌松石樽"
$ws.Range("D23").Value = "Wine Vessel (Zun) with Gold-inlaid Turquoise "
$ws.Range("C23").Value = "Cuo Jin Qian Song Shi Zun"
$ws.Range("E23").Value = "Unkown"
$ws.Range("F23").Value = "The Palace Museum"
$ws.Range("G23").Value = -330
$ws.Range("H23").Value = -221
$ws.Range("J23").Value = "https://www.dpm.org.cn/collection/bronze/229985.html"
$ws.Range("I23").Value = "static/images/错金嵌松石樽.png"

# Match the saved selection state left behind in the worksheet.
$null = $ws.Range("D28").Select()
